$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 2786
$ws.Cells.Item(34, 9).Value = 1631.8
$ws.Cells.Item(34, 10).Value = 6633.3335
$ws.Cells.Item(34, 11).Value = 1631.8
$ws.Cells.Item(34, 12).Value = 6633.3335
$ws.Cells.Item(34, 13).Value = -1428.8
$ws.Cells.Item(34, 14).Value = -7039.3335
$ws.Cells.Item(36, 8).Value = 2786
$ws.Cells.Item(36, 9).Value = 1631.8
$ws.Cells.Item(36, 10).Value = 6633.3335
$ws.Cells.Item(36, 11).Value = 1631.8
$ws.Cells.Item(36, 12).Value = 6633.3335
$ws.Cells.Item(36, 13).Value = -916.8
$ws.Cells.Item(36, 14).Value = -8063.3335
$ws.Cells.Item(64, 8).Value = 6960
$ws.Cells.Item(64, 10).Value = 6960
$ws.Cells.Item(64, 12).Value = 6960
$ws.Cells.Item(64, 14).Value = -7456
$ws.Cells.Item(67, 8).Value = 6960
$ws.Cells.Item(67, 10).Value = 6960
$ws.Cells.Item(67, 12).Value = 6960
$ws.Cells.Item(67, 14).Value = -8676
$ws.Cells.Item(70, 8).Value = 2078.1667
$ws.Cells.Item(70, 9).Value = 2082.1667
$ws.Cells.Item(70, 11).Value = 6246.500100000001
$ws.Cells.Item(70, 13).Value = -5976.500100000001
$ws.Cells.Item(73, 8).Value = 2078.1667
$ws.Cells.Item(73, 9).Value = 2082.1667
$ws.Cells.Item(73, 11).Value = 6246.500100000001
$ws.Cells.Item(73, 13).Value = -5310.500100000001
$ws.Cells.Item(76, 8).Value = 4083.875
$ws.Cells.Item(76, 9).Value = 3385.3076
$ws.Cells.Item(76, 10).Value = 4909.4546
$ws.Cells.Item(76, 11).Value = 3385.3076
$ws.Cells.Item(76, 12).Value = 4909.4546
$ws.Cells.Item(76, 13).Value = -3070.3076
$ws.Cells.Item(76, 14).Value = -5539.4546
$ws.Cells.Item(79, 8).Value = 4083.875
$ws.Cells.Item(79, 9).Value = 3385.3076
$ws.Cells.Item(79, 10).Value = 4909.4546
$ws.Cells.Item(79, 11).Value = 3385.3076
$ws.Cells.Item(79, 12).Value = 4909.4546
$ws.Cells.Item(79, 13).Value = -2293.3076
$ws.Cells.Item(79, 14).Value = -7093.4546
$ws.Cells.Item(92, 8).Value = 307.92856
$ws.Cells.Item(92, 9).Value = 220.63637
$ws.Cells.Item(92, 10).Value = 628
$ws.Cells.Item(92, 11).Value = 220.63637
$ws.Cells.Item(92, 12).Value = 628
$ws.Cells.Item(92, 13).Value = 1027.36363
$ws.Cells.Item(92, 14).Value = -3124
$ws.Cells.Item(107, 8).Value = 4015.8635
$ws.Cells.Item(107, 10).Value = 3237.5
$ws.Cells.Item(107, 12).Value = 3237.5
$ws.Cells.Item(107, 14).Value = -7077.5
$ws.Cells.Item(132, 8).Value = 7726.8
$ws.Cells.Item(132, 9).Value = 7447.864
$ws.Cells.Item(132, 11).Value = 22343.592
$ws.Cells.Item(132, 13).Value = -19813.592
$ws.Cells.Item(137, 8).Value = 11630606
$ws.Cells.Item(137, 9).Value = 13890458
$ws.Cells.Item(137, 10).Value = 8513
$ws.Cells.Item(137, 11).Value = 41671374
$ws.Cells.Item(137, 12).Value = 25539
$ws.Cells.Item(137, 13).Value = -41668824
$ws.Cells.Item(137, 14).Value = -30639
$ws.Cells.Item(138, 8).Value = 2477
$ws.Cells.Item(138, 9).Value = 1350.5
$ws.Cells.Item(138, 10).Value = 2996.923
$ws.Cells.Item(138, 11).Value = 4051.5
$ws.Cells.Item(138, 12).Value = 8990.769
$ws.Cells.Item(138, 13).Value = 1088.5
$ws.Cells.Item(138, 14).Value = -19270.769

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 4669.3335
$ws.Cells.Item(19, 9).Value = 4669.3335
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 4669.3335
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -4440.3335
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(74, 8).Value = 1430262.2
$ws.Cells.Item(74, 9).Value = 1988599.5
$ws.Cells.Item(74, 10).Value = 9040
$ws.Cells.Item(74, 11).Value = 1988599.5
$ws.Cells.Item(74, 12).Value = 9040
$ws.Cells.Item(74, 13).Value = -1987725.5
$ws.Cells.Item(74, 14).Value = -10788
$ws.Cells.Item(77, 8).Value = 1430262.2
$ws.Cells.Item(77, 9).Value = 1988599.5
$ws.Cells.Item(77, 10).Value = 9040
$ws.Cells.Item(77, 11).Value = 9942997.5
$ws.Cells.Item(77, 12).Value = 45200
$ws.Cells.Item(77, 13).Value = -9938629.5
$ws.Cells.Item(77, 14).Value = -53936

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 37212
$ws.Cells.Item(20, 9).Value = 47385.695
$ws.Cells.Item(20, 10).Value = 3784.1428
$ws.Cells.Item(20, 11).Value = 47385.695
$ws.Cells.Item(20, 12).Value = 3784.1428
$ws.Cells.Item(20, 13).Value = -47138.695
$ws.Cells.Item(20, 14).Value = -4278.1428
$ws.Cells.Item(107, 8).Value = 1952.7222
$ws.Cells.Item(107, 9).Value = 1945.4667
$ws.Cells.Item(107, 11).Value = 1945.4667
$ws.Cells.Item(107, 13).Value = -25.46669999999995
$ws.Cells.Item(134, 8).Value = 2733561.8
$ws.Cells.Item(134, 9).Value = 1343.25
$ws.Cells.Item(134, 11).Value = 4029.75
$ws.Cells.Item(134, 13).Value = -1494.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1391824.1
$ws.Cells.Item(31, 9).Value = 1504363.9
$ws.Cells.Item(31, 11).Value = 1504363.9
$ws.Cells.Item(31, 13).Value = -1504068.9
$ws.Cells.Item(34, 8).Value = 1391824.1
$ws.Cells.Item(34, 9).Value = 1504363.9
$ws.Cells.Item(34, 11).Value = 1504363.9
$ws.Cells.Item(34, 13).Value = -1504161.9
$ws.Cells.Item(107, 8).Value = 1190.4445
$ws.Cells.Item(107, 9).Value = 1256.4375
$ws.Cells.Item(107, 10).Value = 662.5
$ws.Cells.Item(107, 11).Value = 1256.4375
$ws.Cells.Item(107, 12).Value = 662.5
$ws.Cells.Item(107, 13).Value = 663.5625
$ws.Cells.Item(107, 14).Value = -4502.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 610.875
$ws.Cells.Item(12, 10).Value = 379.44446
$ws.Cells.Item(12, 12).Value = 1138.33338
$ws.Cells.Item(12, 14).Value = -1484.33338
$ws.Cells.Item(132, 8).Value = 3205.9167
$ws.Cells.Item(132, 9).Value = 3998.75
$ws.Cells.Item(132, 10).Value = 2809.5
$ws.Cells.Item(132, 11).Value = 35988.75
$ws.Cells.Item(132, 12).Value = 25285.5
$ws.Cells.Item(132, 13).Value = -33458.75
$ws.Cells.Item(132, 14).Value = -30345.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 25371.75
$ws.Cells.Item(49, 10).Value = 25706.715
$ws.Cells.Item(49, 12).Value = 25706.715
$ws.Cells.Item(49, 14).Value = -26074.715
$ws.Cells.Item(122, 8).Value = 2339.24
$ws.Cells.Item(122, 9).Value = 1886.375
$ws.Cells.Item(122, 10).Value = 3144.3333
$ws.Cells.Item(122, 11).Value = 5659.125
$ws.Cells.Item(122, 12).Value = 9432.999899999999
$ws.Cells.Item(122, 13).Value = -3209.125
$ws.Cells.Item(122, 14).Value = -14332.9999
$ws.Cells.Item(140, 8).Value = 83333
$ws.Cells.Item(140, 10).Value = 83333
$ws.Cells.Item(140, 12).Value = 83333
$ws.Cells.Item(140, 14).Value = -93693

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3416.2
$ws.Cells.Item(22, 9).Value = 2433.3333
$ws.Cells.Item(22, 10).Value = 3661.9167
$ws.Cells.Item(22, 11).Value = 2433.3333
$ws.Cells.Item(22, 12).Value = 3661.9167
$ws.Cells.Item(22, 13).Value = -2138.3333
$ws.Cells.Item(22, 14).Value = -4251.9167
$ws.Cells.Item(27, 8).Value = 3416.2
$ws.Cells.Item(27, 9).Value = 2433.3333
$ws.Cells.Item(27, 10).Value = 3661.9167
$ws.Cells.Item(27, 11).Value = 2433.3333
$ws.Cells.Item(27, 12).Value = 3661.9167
$ws.Cells.Item(27, 13).Value = -2326.3333
$ws.Cells.Item(27, 14).Value = -3875.9167
$ws.Cells.Item(46, 8).Value = 3033.8462
$ws.Cells.Item(46, 9).Value = 905.6667
$ws.Cells.Item(46, 10).Value = 4363.9585
$ws.Cells.Item(46, 11).Value = 905.6667
$ws.Cells.Item(46, 12).Value = 4363.9585
$ws.Cells.Item(46, 13).Value = -717.6667
$ws.Cells.Item(46, 14).Value = -4739.9585
$ws.Cells.Item(55, 8).Value = 1148.5588
$ws.Cells.Item(55, 9).Value = 1133
$ws.Cells.Item(55, 11).Value = 1133
$ws.Cells.Item(55, 13).Value = -960
$ws.Cells.Item(93, 8).Value = 2489.5217
$ws.Cells.Item(93, 9).Value = 1662.3529
$ws.Cells.Item(93, 10).Value = 4833.1665
$ws.Cells.Item(93, 11).Value = 1662.3529
$ws.Cells.Item(93, 12).Value = 4833.1665
$ws.Cells.Item(93, 13).Value = -414.3529000000001
$ws.Cells.Item(93, 14).Value = -7329.1665

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 5177.5356
$ws.Cells.Item(14, 10).Value = 5939.8823
$ws.Cells.Item(14, 12).Value = 5939.8823
$ws.Cells.Item(14, 14).Value = -6275.8823
$ws.Cells.Item(113, 8).Value = 670.08
$ws.Cells.Item(113, 9).Value = 559.2273
$ws.Cells.Item(113, 11).Value = 1677.6819
$ws.Cells.Item(113, 13).Value = 492.3181
$ws.Cells.Item(132, 8).Value = 3789290
$ws.Cells.Item(132, 9).Value = 4505454.5
$ws.Cells.Item(132, 11).Value = 13516363.5
$ws.Cells.Item(132, 13).Value = -13513833.5
$ws.Cells.Item(140, 8).Value = 79974.5
$ws.Cells.Item(140, 10).Value = 79974.5
$ws.Cells.Item(140, 12).Value = 79974.5
$ws.Cells.Item(140, 14).Value = -90334.5
$ws.Cells.Item(141, 8).Value = 97500
$ws.Cells.Item(141, 10).Value = 97500
$ws.Cells.Item(141, 12).Value = 97500
$ws.Cells.Item(141, 14).Value = -107860
